$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): swap A1/B1 text (name <-> label) ---
$ws.Range("A1").Value = "label"
$ws.Range("B1").Value = "name"

# --- Data rows: update the name (col A) / label (col B) text values ---
$ws.Range("A2").Value = "arbeidsvolume"
$ws.Range("B2").Value = "a__m_mn_ "

$ws.Range("A3").Value = "bbp"
$ws.Range("B3").Value = "bbp_m_wn"

$ws.Range("A4").Value = "belasting"
$ws.Range("B4").Value = "bet_c_wn"

$ws.Range("A5").Value = "coll.lasten"
$ws.Range("B5").Value = "clt_c_wn"

# --- Row 6: the whole row's data is removed (labels + values cleared) ---
$ws.Range("A6:F6").ClearContents()

# --- Column widths: column A gets its own (wider) width; column B keeps its existing width untouched ---
$ws.Columns.Item(1).ColumnWidth = 14.584

# --- Selection moves to A9 ---
$ws.Range("A9").Select()
